# Add the "2022-Q3" fund-holding sheet and update the "总计" (Total) summary
# sheet to include the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet right after "总计" (i.e. before "2022-Q2"),
#    and name it "2022-Q3".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $total)
$newSheet.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2. Seed the new sheet's layout/formatting by copying the header row +
#    one formatted data row from the existing "2022-Q2" sheet (same
#    column layout: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#    持有市值(亿元)/仓位排名). Column A's header cell is intentionally
#    left blank (matches the source sheets), so B:H and A are copied
#    separately.
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2022-Q2")
$dst = $wb.Worksheets.Item("2022-Q3")
$src.Range("B1:H3").Copy($dst.Range("B1:H3"))
$dst = $wb.Worksheets.Item("2022-Q3")
$src.Range("A2:A3").Copy($dst.Range("A2:A3"))

# Extend the styled-row template down to rows 4-7 (6 data rows total).
$dst = $wb.Worksheets.Item("2022-Q3")
$dst.Range("A2:H2").Copy($dst.Range("A4:H4"))
$dst = $wb.Worksheets.Item("2022-Q3")
$dst.Range("A2:H2").Copy($dst.Range("A5:H5"))
$dst = $wb.Worksheets.Item("2022-Q3")
$dst.Range("A2:H2").Copy($dst.Range("A6:H6"))
$dst = $wb.Worksheets.Item("2022-Q3")
$dst.Range("A2:H2").Copy($dst.Range("A7:H7"))

# ---------------------------------------------------------------------
# 3. Fill in the real 2022-Q3 data.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2022-Q3")

# Header row stays the same as the copied template (基金代码, 基金名称,
# 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名) - no change
# needed there.

function Set-TextCell($sheet, $addr, $text) {
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $text
}

# Row 2
$ws.Range("A2").Value = 0
Set-TextCell $ws "B2" "014915"
Set-TextCell $ws "C2" "财通匠心优选一年持有期混合A"
Set-TextCell $ws "D2" "5.65"
Set-TextCell $ws "E2" "81.89"
Set-TextCell $ws "F2" "4.72"
Set-TextCell $ws "G2" "0.2667"
$ws.Range("H2").Value = 9

# Row 3
$ws.Range("A3").Value = 1
Set-TextCell $ws "B3" "011481"
Set-TextCell $ws "C3" "广发瑞锦一年定开混合"
Set-TextCell $ws "D3" "2.62"
Set-TextCell $ws "E3" "89.31"
Set-TextCell $ws "F3" "5.36"
Set-TextCell $ws "G3" "0.1404"
$ws.Range("H3").Value = 7

# Row 4
$ws.Range("A4").Value = 2
Set-TextCell $ws "B4" "501046"
Set-TextCell $ws "C4" "财通多策略福鑫定期开放灵活配置混合"
Set-TextCell $ws "D4" "2.82"
Set-TextCell $ws "E4" "85.55"
Set-TextCell $ws "F4" "4.72"
Set-TextCell $ws "G4" "0.1331"
$ws.Range("H4").Value = 9

# Row 5
$ws.Range("A5").Value = 3
Set-TextCell $ws "B5" "009062"
Set-TextCell $ws "C5" "财通智慧成长混合A"
Set-TextCell $ws "D5" "2.17"
Set-TextCell $ws "E5" "84.78"
Set-TextCell $ws "F5" "4.97"
Set-TextCell $ws "G5" "0.1078"
$ws.Range("H5").Value = 10

# Row 6
$ws.Range("A6").Value = 4
Set-TextCell $ws "B6" "009063"
Set-TextCell $ws "C6" "财通智慧成长混合C"
Set-TextCell $ws "D6" "1.50"
Set-TextCell $ws "E6" "84.78"
Set-TextCell $ws "F6" "4.97"
Set-TextCell $ws "G6" "0.0746"
$ws.Range("H6").Value = 10

# Row 7
$ws.Range("A7").Value = 5
Set-TextCell $ws "B7" "014916"
Set-TextCell $ws "C7" "财通匠心优选一年持有期混合C"
Set-TextCell $ws "D7" "0.61"
Set-TextCell $ws "E7" "81.89"
Set-TextCell $ws "F7" "4.72"
Set-TextCell $ws "G7" "0.0288"
$ws.Range("H7").Value = 9

# ---------------------------------------------------------------------
# 4. Update the "总计" (Total) sheet: insert the 2022-Q3 summary row,
#    pushing 2022-Q2 / 2021-Q3 / 2021-Q2 down by one row.
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Give row 5 the same formatting as the existing data rows before
# writing into it (copies style + border from row 2's A column, etc.).
$tot.Range("A2:D2").Copy($tot.Range("A5:D5"))

# Rewrite rows 2-5 bottom-up so we don't clobber data we still need to
# shift down.
$tot = $wb.Worksheets.Item("总计")
$tot.Range("A5").Value = 3
$tot.Range("B5").Value = "2021-Q2"
$tot.Range("C5").Value = 2
$tot.Range("D5").Value = 0

$tot.Range("A4").Value = 2
$tot.Range("B4").Value = "2021-Q3"
$tot.Range("C4").Value = 3
$tot.Range("D4").Value = 0.24

$tot.Range("A3").Value = 1
$tot.Range("B3").Value = "2022-Q2"
$tot.Range("C3").Value = 2
$tot.Range("D3").Value = 0.11

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q3"
$tot.Range("C2").Value = 6
$tot.Range("D2").Value = 0.75

# ---------------------------------------------------------------------
# 5. Restore the originally-active sheet ("2021-Q2") so the saved
#    workbook's selected tab doesn't change because of the insert.
# ---------------------------------------------------------------------
$orig = $wb.Worksheets.Item("2021-Q2")
$orig.Activate()
